{"js": "// Merge the two \"where ... >= 2012/08/01\" / \"and ... >= 2012/08/31\" paragraphs\n// (XQuery date-range filter) into a single paragraph using a quantified\n// \"some ... satisfies contains(...)\" expression, per the target diff:\n//\n//   where $entidad/actosfesteros/acto/fecha >= \"2012/08/01\"\n//   and $entidad/actosfesteros/acto/fecha >= \"2012/08/31\"\n// becomes\n//   where some $x in $entidad/actosfesteros/acto/fecha satisfies contains($x, \"/8/\")\n//\n// Locate the two paragraphs robustly by searching for the unique literal\n// text they contain (rather than relying on a hard-coded paragraph index).\nconst body = context.document.body;\n\nconst startHits = body.search('>= \"2012/08/01\"', { matchCase: true });\nstartHits.load(\"items\");\nawait context.sync();\n\nif (startHits.items.length === 0) {\n  throw new Error('Could not find the paragraph containing \\'>= \"2012/08/01\"\\'');\n}\n\nconst firstPara = startHits.items[0].paragraphs.getFirst();\nfirstPara.load(\"text\");\nawait context.sync();\n\nconst secondPara = firstPara.getNext();\nsecondPara.load(\"text\");\nawait context.sync();\n\n// Sanity-check we grabbed the expected two paragraphs before rewriting them.\nif (!/2012\\/08\\/01/.test(firstPara.text) || !/2012\\/08\\/31/.test(secondPara.text)) {\n  throw new Error(\"Unexpected paragraph contents while locating the C8 where/and clauses\");\n}\n\n// Range spanning both paragraphs (start of the first through end of the\n// second) so the replacement collapses them into a single paragraph.\nconst rangeStart = firstPara.getRange(Word.RangeLocation.start);\nconst rangeEnd = secondPara.getRange(Word.RangeLocation.end);\nconst combinedRange = rangeStart.expandTo(rangeEnd);\n\n// Build the replacement as a literal OOXML fragment so the exact run /\n// proofErr / bookmark layout from the target revision is reproduced\n// (the existing \"_GoBack\" bookmark is preserved in place).\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>where</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>some</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> $x in</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> $entidad/</w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>actosfesteros</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\">/acto/fecha </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>satisfies</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>contains</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t>($x,</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:r><w:t>&quot;</w:t></w:r>\n            <w:r><w:t>/8/</w:t></w:r>\n            <w:r><w:t>&quot;</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ncombinedRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Merge the two \"where ... >= 2012/08/01\" / \"and ... >= 2012/08/31\" paragraphs\n# (XQuery date-range filter) into a single paragraph using a quantified\n# \"some ... satisfies contains(...)\" expression, per the target diff:\n#\n#   where $entidad/actosfesteros/acto/fecha >= \"2012/08/01\"\n#   and $entidad/actosfesteros/acto/fecha >= \"2012/08/31\"\n# becomes\n#   where some $x in $entidad/actosfesteros/acto/fecha satisfies contains($x, \"/8/\")\n#\n$d = $word.ActiveDocument\n\n# Locate the two paragraphs robustly via Find rather than a hard-coded index.\n$finder = $d.Content\n$found = $finder.Find.Execute('>= \"2012/08/01\"')\nif (-not $found) {\n    throw 'Could not find the paragraph containing >= \"2012/08/01\"'\n}\n\n$para1 = $finder.Paragraphs(1)\n$para2 = $para1.Next()\n\nif ($para1.Range.Text -notlike '*2012/08/01*' -or $para2.Range.Text -notlike '*2012/08/31*') {\n    throw 'Unexpected paragraph contents while locating the C8 where/and clauses'\n}\n\n# Range spanning both paragraphs (start of the first through end of the\n# second, including the paragraph mark) so the replacement collapses them\n# into a single paragraph. Re-fetch a fresh Range via $d.Range(...) (rather\n# than reusing/duplicating the Find range) so InsertXML replaces the span\n# instead of inserting after it.\n$fullRange = $d.Range($para1.Range.Start, $para2.Range.End)\n\n# Build the replacement as a literal OOXML fragment so the exact run /\n# proofErr / bookmark layout from the target revision is reproduced\n# (the existing \"_GoBack\" bookmark is preserved in place).\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>where</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>some</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> `$x in</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> `$entidad/</w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>actosfesteros</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\">/acto/fecha </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>satisfies</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>contains</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t>(`$x,</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:r><w:t>&quot;</w:t></w:r>\n            <w:r><w:t>/8/</w:t></w:r>\n            <w:r><w:t>&quot;</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$fullRange.InsertXML($ooxml)\n"}
